$d = $word.ActiveDocument

$pairs = @(
    @("128×3=384", "977×3=2931"),
    @("452×7=3164", "689×2=1378"),
    @("516×7=3612", "651×9=5859"),
    @("646×8=5168", "806×7=5642"),
    @("527×3=1581", "696×7=4872"),
    @("337×8=2696", "497×6=2982"),
    @("813×4=3252", "305×5=1525"),
    @("139×5=695", "671×4=2684"),
    @("845×4=3380", "266×2=532"),
    @("440×9=3960", "169×5=845"),
    @("488×9=4392", "606×4=2424"),
    @("342×3=1026", "593×7=4151"),
    @("121×5=605", "844×2=1688"),
    @("867×2=1734", "712×5=3560"),
    @("548×2=1096", "290×5=1450"),
    @("867×4=3468", "974×2=1948"),
    @("409×8=3272", "121×3=363"),
    @("186×7=1302", "189×2=378"),
    @("747×3=2241", "713×6=4278"),
    @("694×5=3470", "605×6=3630"),
    @("658×7=4606", "646×7=4522"),
    @("873×2=1746", "305×9=2745"),
    @("745×6=4470", "435×6=2610"),
    @("942×8=7536", "773×2=1546"),
    @("303×2=606", "806×5=4030")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
